$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend / update the "kommentar" (column C) notes on several rows ---
# (written in the same relative order the shared-string table ends up with
# in the target file: skript, lars, fatt)

# Row 11: "GIS-data lyornas avstånd andra fjällrävslyor" — status moves from
# "ej påbörjat" to "påbörjat" (copy the "klar" row's special font first so we
# keep the "Calibri (Brödtext)" font face, then recolor it), and the comment
# gets a big extension.
$ws.Range("B9").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Value = "påbörjat"
$ws.Range("B11").Font.ThemeColor = 9
$ws.Range("C11").Value = 'Fått skript av Rasmus, ändra om till mina data. Fattas för vissa år i Rasmus.kulldata till Tor och i Lypositioner kullar 2000-2017 SWEREF99 per kull (finns mer data i den senare även fast Rasmus tagit bort lågår och omatade lyor. Komplettera med data från filen Red fox feeding och hitta vilka två lyor det var reproduktion på 2009.'

# Row 7: "lyaktivitet vår, Lars data" comment — append follow-up note.
$ws.Range("C7").Value = 'Lars skulle se vad han hade. Han har mailat en del data tidigare år (c:a 2000 - 2005)  till Tomas Meijer och kanske till Anders. Gick inte att få ut något vettigt ur databasen. Det fattas årtal och datum på majoriteten av lybesöken.'

# Row 8: "lyaktivitet reproduktion" comment — append follow-up note.
$ws.Range("C8").Value = 'fått datan av Rasmus. Gjort en separat fil för de ripinventerade lyorna, dock verkar det saknas data för "2003", "2006", "2009", "2012", "2016". 2016 tog jag ur rovbase. 2012 var det inga. 2009 var det två. 2003 och 2006 finns i Peters fil (inlagt). Mailat Alva för att få ett utdrag ur databasen istället.'

$excel.CutCopyMode = $false

# --- Worksheet view: selection cell moved from C20 to C15 ---
$ws.Range("C15").Select() | Out-Null

# --- Page setup: paper size A4 / portrait orientation now specified ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
